# Apply the "HISTORIA DE USUARIOA.docx" edits:
#  - fix a handful of missing Spanish accents
#  - remove the stray "_GoBack" bookmark left over from a previous save
#
# wdReplaceAll = 2 (replace every match found by Find.Execute)
$wdReplaceAll = 2
$wdFindContinue = 1

$d = $word.ActiveDocument

function Fix-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, `
        $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# 1) "Para agilizar la busqueda del articulo." -> "... búsqueda del artículo."
Fix-Text "Para agilizar la busqueda del articulo." "Para agilizar la búsqueda del artículo."

# 2) "La categoria bicicletas" -> "La categoría bicicletas"
Fix-Text "La categoria bicicletas" "La categoría bicicletas"

# 3) "... correspondientes a la categoria seleccionada." -> "... categoría seleccionada."
Fix-Text "categoria seleccionada." "categoría seleccionada."

# 4) "agregar y quitar articulos" -> "... artículos"
Fix-Text "agregar y quitar articulos" "agregar y quitar artículos"

# 5) "tarjetas de credito y debito." -> "... crédito y débito."
Fix-Text "de credito y debito." "de crédito y débito."

# 6) Remove the leftover "_GoBack" bookmark (Word re-creates this automatically at the
#    last edit point; it is stray, unneeded metadata and was removed from the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
